$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 8 (pushes old rows 8-18 down to 10-20)
$ws.Rows("8:9").Insert()

# Populate the new row 8 with the new asserted_distribution test-case data
$ws.Cells.Item(8, 1).Value = "Ursus arctos horribilis"
$ws.Cells.Item(8, 2).Value = 45678
$ws.Cells.Item(8, 3).Value = "California"

# Match the row height / selection state captured in the saved workbook
$ws.Rows(8).RowHeight = 45

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("I8").Select()
